$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header typo
$ws.Range("G1").Value = "¿De Ñuble?"

# Update row 2 with new values
$ws.Range("A2").Value = "test"
$ws.Range("B2").Value = "Nitrógeno Gaseoso"
$ws.Range("C2").Value = "10m3"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2023-11-15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2023-11-15"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "Si"

# Remove rows 3 and 4 (old history entries)
$ws.Range("A3:G4").Delete()
